# tube.xlsx - jakeR, doncina01, doncina02 tube additions
# (also renames existing "jake" tube set to "jake_B")

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Tube")

# --- Rename existing jake_* rows (21-23) to jake_B_* ---
# Row 21: jake style -> jake_B style (also damage 10 -> 20)
$ws.Range("B21").Value = "jake_B_style"
$ws.Range("C21").Value = "제이크B 스타일"
$ws.Range("G21").Value = "jake_B_attack"
$ws.Range("J21").Value = 20

# Row 22: jake enhancer -> jake_B enhancer
$ws.Range("B22").Value = "jake_B_enhancer"
$ws.Range("C22").Value = "제이크B 인핸서"
$ws.Range("M22").Value = "{(jake_B_attack, 4, 4)}"

# Row 23: jake cooler -> jake_B cooler
$ws.Range("B23").Value = "jake_B_cooler"
$ws.Range("C23").Value = "제이크B 쿨러"

# --- New rows: jake_R tube set (24-26) ---
$ws.Range("A24").Value = 5109
$ws.Range("B24").Value = "jake_R_style"
$ws.Range("C24").Value = "제이크R 스타일"
$ws.Range("D24").Value = "style"
$ws.Range("E24").Value = "Weakness"
$ws.Range("F24").Value = "A"
$ws.Range("G24").Value = "jake_R_attack"
$ws.Range("H24").Value = "{(0.3)}"
$ws.Range("I24").Value = "range"
$ws.Range("J24").Value = 10
$ws.Range("L24").Value = "{(jake_R_skill, 0,5)}"

$ws.Range("A25").Value = 5110
$ws.Range("B25").Value = "jake_R_enhancer"
$ws.Range("C25").Value = "제이크R 인핸서"
$ws.Range("D25").Value = "enhancer"
$ws.Range("E25").Value = "Weakness"
$ws.Range("F25").Value = "A"
$ws.Range("O25").Value = 0

$ws.Range("A26").Value = 5111
$ws.Range("B26").Value = "jake_R_cooler"
$ws.Range("C26").Value = "제이크R 쿨러"
$ws.Range("D26").Value = "cooler"
$ws.Range("E26").Value = "Weakness"
$ws.Range("F26").Value = "A"
$ws.Range("P26").Value = 2

# --- New rows: doncina01 tube set (27-29) ---
$ws.Range("A27").Value = 5112
$ws.Range("B27").Value = "doncina01_style"
$ws.Range("C27").Value = "돈시나01 스타일"
$ws.Range("D27").Value = "style"
$ws.Range("E27").Value = "gangster"
$ws.Range("F27").Value = "A"
$ws.Range("G27").Value = "doncina_skill00"
$ws.Range("H27").Value = "none"
$ws.Range("I27").Value = "melee"
$ws.Range("J27").Value = 50

$ws.Range("A28").Value = 5113
$ws.Range("B28").Value = "doncina01_enhancer"
$ws.Range("C28").Value = "돈시나01 인핸서"
$ws.Range("D28").Value = "enhancer"
$ws.Range("E28").Value = "gangster"
$ws.Range("F28").Value = "B"
$ws.Range("K28").Value = "{(doncinal, 0, 1)}"
$ws.Range("O28").Value = 0

$ws.Range("A29").Value = 5114
$ws.Range("B29").Value = "doncina01_cooler"
$ws.Range("C29").Value = "돈시나01 쿨러"
$ws.Range("D29").Value = "cooler"
$ws.Range("E29").Value = "gangster"
$ws.Range("F29").Value = "B"
$ws.Range("P29").Value = 5

# --- New rows: doncina02 tube set (30-33) ---
$ws.Range("A30").Value = 5115
$ws.Range("B30").Value = "doncina02_style"
$ws.Range("C30").Value = "돈시나02 스타일"
$ws.Range("D30").Value = "style"
$ws.Range("E30").Value = "gangster"
$ws.Range("F30").Value = "B"
$ws.Range("G30").Value = "doncina_skill01"
$ws.Range("H30").Value = "{(1.0)}"
$ws.Range("I30").Value = "range"
$ws.Range("J30").Value = 10

$ws.Range("A31").Value = 5116
$ws.Range("B31").Value = "doncina02_enhancer"
$ws.Range("C31").Value = "돈시나02 인핸서"
$ws.Range("D31").Value = "enhancer"
$ws.Range("E31").Value = "gangster"
$ws.Range("F31").Value = "A"
$ws.Range("L31").Value = "{(doncinal, 1, 5)}"
$ws.Range("O31").Value = 0

$ws.Range("A32").Value = 5117
$ws.Range("B32").Value = "doncina02_cooler"
$ws.Range("C32").Value = "돈시나02 쿨러"
$ws.Range("D32").Value = "cooler"
$ws.Range("E32").Value = "gangster"
$ws.Range("F32").Value = "B"
$ws.Range("P32").Value = 10

$ws.Range("A33").Value = 5118
$ws.Range("B33").Value = "doncina02_relic"
$ws.Range("C33").Value = "돈시나02 렐릭"
$ws.Range("D33").Value = "relic"
$ws.Range("E33").Value = "gangster"
$ws.Range("F33").Value = "C"
$ws.Range("Q33").Value = "snare_01"
